$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 14, shifting the existing rows (14-22) down to (15-23).
# Using EntireRow.Insert() copies the formatting from the row above (row 13),
# matching Excel's default "insert row" behavior.
$ws.Rows(14).Insert()

# Fill in the new row 14 with the meeting follow-up task.
$ws.Range("A14").Value = [DateTime]::Parse("2017-03-08")
$ws.Range("B14").Value = [DateTime]::Parse("12:00:00")
$ws.Range("C14").Value = "Group Meeting"
$ws.Range("D14").Value = "Finish initial research, draw schematics, buy parts"
$ws.Range("E14").Value = ""

# Update the view so the new row is visible/selected, matching the saved view state.
$ws.Range("D14").Select()
$excel.ActiveWindow.ScrollRow = 8
